$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 256, shifting existing rows 256:274 down to 257:275.
$ws.Rows.Item(256).Insert()

# Populate the newly inserted row 256 with the new weekly record.
$ws.Range("A256").Value = 7
$ws.Range("B256").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C256").Value = "Ñuble"
$ws.Range("D256").Value = 44931
$ws.Range("E256").Value = 16
$ws.Range("F256").Value = 100112017
$ws.Range("G256").Value = "Apio"
$ws.Range("H256").Value = "Americana (o)"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 80
$ws.Range("K256").Value = 10000
$ws.Range("L256").Value = 10000
$ws.Range("M256").Value = 10000
$ws.Range("N256").Value = "$/docena de matas"
$ws.Range("O256").Value = "Provincia del Elquí"
$ws.Range("P256").Value = 1667
$ws.Range("Q256").Value = 6
$ws.Range("R256").Value = "Hortaliza"
